$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.758.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.109.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.16%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "388.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.584.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.100.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.983"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.893.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0450"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.79%  "
$ws.Range("E39").Value = "  +6.78%  "
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").Value = "  -4.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.57%  "
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.048.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.408.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("E51").Value = "  +6.00%  "
